# Add the "training milestone" worksheet (commit: "add training milestone in Result.xlsx")
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after "Mask RCNN R100" so it becomes the 2nd
# (active) tab, matching sheetId="2" / activeTab="1" in the target workbook.
$ws = $wb.Worksheets.Add($null, $ws1)
$ws.Name = "training milestone"

# ---- Row 1: title ----
$ws.Range("A1").Value = "command"

# ---- Rows 2-4: the training command, wrapped across three merged rows ----
$ws.Range("A2").Value = 'python tools/train_net.py \'
$ws.Range("A3").Value = '    --cfg configs/wattanapong/train_e2e_mask_rcnn_R-101-FPN_4x_coco2014_train_valminusmini_2gpu.yaml \'
$ws.Range("A4").Value = '    OUTPUT_DIR /tmp/detectron-output > train_test_e2e_mask_rcnn_R-101-FPN_4x_coco2014_train_valminusmini_lr0_02.txt'

$ws.Range("A2:G2").Merge()
$ws.Range("A3:G3").Merge()
$ws.Range("A4:G4").Merge()
$ws.Range("A2:G4").HorizontalAlignment = -4131
$ws.Range("A2:G4").VerticalAlignment = -4108

# ---- Row 7/9/10/11/12/14 detail cells (entered in this order by the author) ----
$ws.Range("C7").Value = "tools"
$ws.Range("D7").Value = "train_net.py"

$ws.Range("C9").Value = "detectron"
$ws.Range("D9").Value = "utils"
$ws.Range("E9").Value = "train.py"
$ws.Range("F9").Value = "train_model()"

$ws.Range("G10").Value = "create_model()"

$ws.Range("E7").Value = "detectron.utils.train.train_model()"

$ws.Range("H11").Value = "checkpoint_iter"

$ws.Range("H12").Value = "model_builder.create"

$ws.Range("D14").Value = "modeling"

# ---- Row 6: table header (typed after the detail rows) ----
$ws.Range("A6").Value = "No."
$ws.Range("B6").Value = "Line"
$ws.Range("C6").Value = "directory & function"

# ---- Row 14 remaining cells ----
$ws.Range("C14").Value = "detectron"
$ws.Range("E14").Value = "model_builder.py"
$ws.Range("F14").Value = "create()"

# ---- Numeric cells ----
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 110
$ws.Range("A9").Value = 2
$ws.Range("B9").Value = 50
$ws.Range("B10").Value = 98
$ws.Range("B11").Value = 119
$ws.Range("B12").Value = 134
$ws.Range("A14").Value = 3
$ws.Range("B14").Value = 107

# ---- Merges ----
$ws.Range("C6:M6").Merge()
$ws.Range("E7:G7").Merge()

# ---- Fills (yellow) ----
$ws.Range("E7:G7").Interior.Color = 65535
$ws.Range("E8:G8").Interior.Color = 65535
$ws.Range("F9").Interior.Color = 65535
$ws.Range("G10").Interior.Color = 65535
$ws.Range("F14").Interior.Color = 65535

# ---- Alignment ----
$ws.Range("C6:M6").HorizontalAlignment = -4108
$ws.Range("E7:G8").HorizontalAlignment = -4108

# ---- Column widths (best-fit-ish, matches the recorded widths) ----
$ws.Columns.Item(1).ColumnWidth = 26.77734375
$ws.Columns.Item(2).ColumnWidth = 12.109375
$ws.Columns.Item(4).ColumnWidth = 10.88671875
$ws.Columns.Item(5).ColumnWidth = 15.21875
$ws.Columns.Item(6).ColumnWidth = 13.5546875
$ws.Columns.Item(7).ColumnWidth = 13.5546875
$ws.Columns.Item(8).ColumnWidth = 18.6640625

$ws.PageSetup.Orientation = 1

# ---- Selection / view state on the new sheet ----
$ws.Range("G19").Select()

Write-Output "training milestone sheet added"
